# Added Configurable zero_before_threshold parameter to enable setting dims
# before noise_threshold or First Rise Point to 0.
# This updates the First_Noticeable_Increase_Index (C), the
# First_Noticeable_Increase_Cumulative_Value (E) and the resulting Pulse_Width (G)
# on each of the Step3_DataPts_* sheets to reflect the new zero_before_threshold
# behavior.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> row -> updated C / E / G values
$updates = @{
    "Step3_DataPts_0.5" = @{
        2 = @{ C = 87; E = 0.008441812073723569; G = 41 }
        3 = @{ C = 87; E = 0.004705503334542409; G = 43 }
        4 = @{ C = 87; E = 0.005337774033308033; G = 43 }
        5 = @{ C = 88; E = 0.006026037585510682; G = 42 }
        6 = @{ C = 87; E = 0.005940632590531252; G = 42 }
    }
    "Step3_DataPts_0.7" = @{
        2 = @{ C = 87; E = 0.008441812073723569; G = 53 }
        3 = @{ C = 87; E = 0.004705503334542409; G = 53 }
        4 = @{ C = 87; E = 0.005337774033308033; G = 54 }
        5 = @{ C = 88; E = 0.006026037585510682; G = 52 }
        6 = @{ C = 87; E = 0.005940632590531252; G = 52 }
    }
    "Step3_DataPts_0.8" = @{
        2 = @{ C = 87; E = 0.008441812073723569; G = 62 }
        3 = @{ C = 87; E = 0.004705503334542409; G = 63 }
        4 = @{ C = 87; E = 0.005337774033308033; G = 63 }
        5 = @{ C = 88; E = 0.006026037585510682; G = 62 }
        6 = @{ C = 87; E = 0.005940632590531252; G = 62 }
    }
    "Step3_DataPts_0.9" = @{
        2 = @{ C = 87; E = 0.008441812073723569; G = 77 }
        3 = @{ C = 87; E = 0.004705503334542409; G = 78 }
        4 = @{ C = 87; E = 0.005337774033308033; G = 77 }
        5 = @{ C = 88; E = 0.006026037585510682; G = 77 }
        6 = @{ C = 87; E = 0.005940632590531252; G = 77 }
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($rowNum in $rows.Keys) {
        $vals = $rows[$rowNum]
        $ws.Cells.Item($rowNum, 3).Value = $vals.C
        $ws.Cells.Item($rowNum, 5).Value = $vals.E
        $ws.Cells.Item($rowNum, 7).Value = $vals.G
    }
}
